$d = $word.ActiveDocument

# --- Remove the stray "_GoBack" bookmark from the title-image paragraph;
#     it will be re-created around the new trailing empty paragraph below.
$d.Bookmarks("_GoBack").Delete()

# --- Replace the trailing empty paragraph with the new title-page / abstract /
#     introduction content, expressed as literal WordprocessingML via InsertXML.
#     InsertXML inserts the supplied paragraphs immediately before the (now
#     collapsed) range, leaving the original empty paragraph as the new last one.
$p = $d.Paragraphs.Last
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:sz w:val="44"/>
          <w:szCs w:val="44"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia"/>
          <w:b/>
          <w:sz w:val="44"/>
          <w:szCs w:val="44"/>
        </w:rPr>
        <w:t>计算机系统基础论文</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:ascii="楷体" w:hAnsi="楷体" w:eastAsia="楷体"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="楷体" w:hAnsi="楷体" w:eastAsia="楷体"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>作者：</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="楷体" w:hAnsi="楷体" w:eastAsia="楷体"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>丁楷轩  学号：20151104711  指导教师：朝力萌</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>摘要：</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>本论文介绍计算机系统的主要组成及简单的工作流程、c语言中的指针的用法和作用和Arduino系统的简单介绍。通过这些来总结计算机系统的主要知识。</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>关键词：计算机系统、计算机软件、指针、Arduino。</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>The computer system base paper</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>Yuan Dongqiang</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>Abstract:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>This paper introduces the main computer system composition and the simple working process, and the usage and function of the pointer in c language and a brief introduction of the Arduino system. To summarize the main computer system by the knowledge.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="仿宋" w:hAnsi="仿宋" w:eastAsia="仿宋"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
        </w:rPr>
        <w:t>Key words: computer systems, computer software, pointer, the Arduino.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>引言：</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:ascii="黑体" w:hAnsi="黑体" w:eastAsia="黑体"/>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>通过半年的学习，让我们认识了计算机系统的主要</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>组成部分及其原理</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>，同时学会了如何使用Arduino，以及通过Arduino制作一些简单的机器人，一些简单的程序控制的模块。通过这半年的学习改变了我们一开始所认识的计算机就是电脑的思想，对计算机的结构也有了简单的认识，工作原理</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>也</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>有了一些</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
          <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
        </w:rPr>
        <w:t>更为深刻的</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="eastAsia" w:asciiTheme="majorEastAsia" w:hAnsiTheme="majorEastAsia" w:eastAsiaTheme="majorEastAsia"/>
          <w:szCs w:val="21"/>
        </w:rPr>
        <w:t>了解。</w:t>
      </w:r>
    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p.Range.InsertXML($xml)

# --- Re-create the "_GoBack" bookmark around the (still empty) final paragraph.
$last = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $last.Range)

# --- Mark the "Default Paragraph Font" character style as a quick/recommended
#     style (adds <w:qFormat/> to its style definition).
$d.Styles("Default Paragraph Font").QuickStyle = $true
